$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of column letters to column indices
function Set-CellText($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# Row 2
Set-CellText 2 4 "60.095.63"
Set-CellText 2 5 "  +2.64%  "

# Row 3
Set-CellText 3 4 "3.202.59"
Set-CellText 3 5 "  +1.45%  "

# Row 5
Set-CellText 5 4 "536.99"
Set-CellText 5 5 "  -0.12%  "

# Row 6
Set-CellText 6 4 "144.88"
Set-CellText 6 5 "  +3.41%  "

# Row 7
Set-CellText 7 4 "1.00"
Set-CellText 7 5 "  +0.05%  "

# Row 8
Set-CellText 8 4 "0.528"
Set-CellText 8 5 "  +3.28%  "

# Row 9
Set-CellText 9 4 "7.32"
Set-CellText 9 5 "  -0.36%  "

# Row 11
Set-CellText 11 4 "0.430"
Set-CellText 11 5 "  +2.32%  "

# Row 12
Set-CellText 12 4 "3.755.67"
Set-CellText 12 5 "  +1.48%  "

# Row 13
Set-CellText 13 5 "  -1.20%  "

# Row 14
Set-CellText 14 4 "25.88"
Set-CellText 14 5 "  +0.33%  "

# Row 15
Set-CellText 15 5 "  +1.64%  "

# Row 16
Set-CellText 16 4 "60.120.71"
Set-CellText 16 5 "  +2.56%  "

# Row 17
Set-CellText 17 4 "3.215.00"
Set-CellText 17 5 "  +1.88%  "

# Row 18
Set-CellText 18 4 "6.24"
Set-CellText 18 5 "  +0.28%  "

# Row 19
Set-CellText 19 4 "13.11"
Set-CellText 19 5 "  +0.62%  "

# Row 20
Set-CellText 20 4 "8.29"
Set-CellText 20 5 "  +0.58%  "

# Row 21
Set-CellText 21 4 "375.82"
Set-CellText 21 5 "  +0.75%  "

# Row 22
Set-CellText 22 4 "0.999"
Set-CellText 22 5 "  -0.12%  "

# Row 23
Set-CellText 23 4 "0.523"
Set-CellText 23 5 "  +1.72%  "

# Row 24
Set-CellText 24 4 "70.01"
Set-CellText 24 5 "  -0.04%  "

# Row 25
Set-CellText 25 5 "  +1.27%  "

# Row 26
Set-CellText 26 5 "  +8.62%  "

# Row 27
Set-CellText 27 5 "  +0.56%  "

# Row 28
Set-CellText 28 5 "  +2.95%  "

# Row 29
Set-CellText 29 5 "  +0.30%  "

# Row 30
Set-CellText 30 4 "22.34"
Set-CellText 30 5 "  +1.84%  "

# Row 31
Set-CellText 31 4 "6.13"
Set-CellText 31 5 "  -0.86%  "

# Row 32
Set-CellText 32 4 "5.37"
Set-CellText 32 5 "  +3.42%  "

# Row 33
Set-CellText 33 2 "Fetch.AI"
Set-CellText 33 3 "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-CellText 33 4 "1.20"
Set-CellText 33 5 "  +2.56%  "

# Row 34
Set-CellText 34 2 "Aptos"
Set-CellText 34 3 "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-CellText 34 4 "6.64"
Set-CellText 34 5 "  +6.48%  "

# Row 35
Set-CellText 35 4 "156.82"
Set-CellText 35 5 "  -2.20%  "

# Row 36
Set-CellText 36 5 "  -1.69%  "

# Row 37
Set-CellText 37 4 "2.798.99"
Set-CellText 37 5 "  +6.09%  "

# Row 38
Set-CellText 38 4 "25.56"
Set-CellText 38 5 "  +1.11%  "

# Row 39
Set-CellText 39 4 "0.0704"
Set-CellText 39 5 "  +3.05%  "

# Row 40
Set-CellText 40 5 "  +0.77%  "

# Row 41
Set-CellText 41 5 "  +1.23%  "

# Row 42
Set-CellText 42 4 "39.87"
Set-CellText 42 5 "  +2.82%  "

# Row 43
Set-CellText 43 4 "0.0293"
Set-CellText 43 5 "  +4.30%  "

# Row 44
Set-CellText 44 5 "  +1.09%  "

# Row 45
Set-CellText 45 5 "  +3.66%  "

# Row 46
Set-CellText 46 4 "3.244.95"
Set-CellText 46 5 "  +1.44%  "

# Row 47
Set-CellText 47 4 "0.984"
Set-CellText 47 5 "  +0.03%  "

# Row 48
Set-CellText 48 4 "0.811"
Set-CellText 48 5 "  +6.94%  "

# Row 49
Set-CellText 49 5 "  -1.03%  "

# Row 50
Set-CellText 50 4 "20.59"
Set-CellText 50 5 "  +1.52%  "

# Row 51
Set-CellText 51 5 "  -0.03%  "
